$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 102, shifting existing rows 102..182 down to 103..183.
$ws.Rows.Item(102).Insert()

# Populate the newly inserted row 102 with the new record's data.
$ws.Cells.Item(102, 1).Value = 10
$ws.Cells.Item(102, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(102, 3).Value = "La Araucanía"
$ws.Cells.Item(102, 4).Value = 44741
$ws.Cells.Item(102, 5).Value = 9
$ws.Cells.Item(102, 6).Value = "Fruta"
$ws.Cells.Item(102, 7).Value = 100104
$ws.Cells.Item(102, 8).Value = "Frutos de pepita"
$ws.Cells.Item(102, 9).Value = 100104003
$ws.Cells.Item(102, 10).Value = "Membrillo"
$ws.Cells.Item(102, 11).Value = "Champion"
$ws.Cells.Item(102, 12).Value = "Primera"
$ws.Cells.Item(102, 13).Value = 2
$ws.Cells.Item(102, 14).Value = 200000
$ws.Cells.Item(102, 15).Value = 200000
$ws.Cells.Item(102, 16).Value = 200000
$ws.Cells.Item(102, 17).Value = "$/bins (450 kilos)"
$ws.Cells.Item(102, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(102, 19).Value = 444
$ws.Cells.Item(102, 20).Value = 450
